$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 67.77251700000001
$ws.Range("H2").Value = 203.317551
$ws.Range("I2").Value = 0.4079637943863715
$ws.Range("J2").Value = 0.4079637943863715
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 87.94215800000001
$ws.Range("N2").Value = 263.826474
$ws.Range("O2").Value = 0.8507690866039653
$ws.Range("P2").Value = 0.8507690866039652
$ws.Range("Q2").Value = 5960.061398071687
$ws.Range("R2").Value = 53640.55258264518
$ws.Range("S2").Value = 0.3470829847175812
$ws.Range("T2").Value = 0.3470829847175811

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 67.77251700000001
$ws.Range("H3").Value = 203.317551
$ws.Range("I3").Value = 0.4079637943863715
$ws.Range("J3").Value = 0.4079637943863715
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 13.642319
$ws.Range("N3").Value = 40.926957
$ws.Range("O3").Value = 0.1319783769098539
$ws.Range("P3").Value = 0.1319783769098539
$ws.Range("Q3").Value = 924.5742963469231
$ws.Range("R3").Value = 8321.168667122307
$ws.Range("S3").Value = 0.05384239942109866
$ws.Range("T3").Value = 0.05384239942109866

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 67.77251700000001
$ws.Range("H4").Value = 203.317551
$ws.Range("I4").Value = 0.4079637943863715
$ws.Range("J4").Value = 0.4079637943863715
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.783357333333333
$ws.Range("N4").Value = 5.350072
$ws.Range("O4").Value = 0.01725253648618087
$ws.Range("P4").Value = 0.01725253648618087
$ws.Range("Q4").Value = 120.862615190408
$ws.Range("R4").Value = 1087.763536713672
$ws.Range("S4").Value = 0.007038410247691665
$ws.Range("T4").Value = 0.007038410247691665

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 60.97760633333333
$ws.Range("H5").Value = 182.932819
$ws.Range("I5").Value = 0.3670611149405164
$ws.Range("J5").Value = 0.3670611149405164
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 87.94215800000001
$ws.Range("N5").Value = 263.826474
$ws.Range("O5").Value = 0.8507690866039653
$ws.Range("P5").Value = 0.8507690866039652
$ws.Range("Q5").Value = 5362.502290627801
$ws.Range("R5").Value = 48262.52061565021
$ws.Range("S5").Value = 0.3122842494857763
$ws.Range("T5").Value = 0.3122842494857762

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 60.97760633333333
$ws.Range("H6").Value = 182.932819
$ws.Range("I6").Value = 0.3670611149405164
$ws.Range("J6").Value = 0.3670611149405164
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 13.642319
$ws.Range("N6").Value = 40.926957
$ws.Range("O6").Value = 0.1319783769098539
$ws.Range("P6").Value = 0.1319783769098539
$ws.Range("Q6").Value = 831.8759574557537
$ws.Range("R6").Value = 7486.883617101783
$ws.Range("S6").Value = 0.04844413017657067
$ws.Range("T6").Value = 0.04844413017657066

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 60.97760633333333
$ws.Range("H7").Value = 182.932819
$ws.Range("I7").Value = 0.3670611149405164
$ws.Range("J7").Value = 0.3670611149405164
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 1.783357333333333
$ws.Range("N7").Value = 5.350072
$ws.Range("O7").Value = 0.01725253648618087
$ws.Range("P7").Value = 0.01725253648618087
$ws.Range("Q7").Value = 108.7448614236631
$ws.Range("R7").Value = 978.703752812968
$ws.Range("S7").Value = 0.00633273527816949
$ws.Range("T7").Value = 0.006332735278169489

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 37.37372866666667
$ws.Range("H8").Value = 112.121186
$ws.Range("I8").Value = 0.2249750906731122
$ws.Range("J8").Value = 0.2249750906731122
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 87.94215800000001
$ws.Range("N8").Value = 263.826474
$ws.Range("O8").Value = 0.8507690866039653
$ws.Range("P8").Value = 0.8507690866039652
$ws.Range("Q8").Value = 3286.72635145313
$ws.Range("R8").Value = 29580.53716307817
$ws.Range("S8").Value = 0.1914018524006079
$ws.Range("T8").Value = 0.1914018524006079

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 37.37372866666667
$ws.Range("H9").Value = 112.121186
$ws.Range("I9").Value = 0.2249750906731122
$ws.Range("J9").Value = 0.2249750906731122
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 13.642319
$ws.Range("N9").Value = 40.926957
$ws.Range("O9").Value = 0.1319783769098539
$ws.Range("P9").Value = 0.1319783769098539
$ws.Range("Q9").Value = 509.8643286901114
$ws.Range("R9").Value = 4588.778958211003
$ws.Range("S9").Value = 0.02969184731218455
$ws.Range("T9").Value = 0.02969184731218455

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 37.37372866666667
$ws.Range("H10").Value = 112.121186
$ws.Range("I10").Value = 0.2249750906731122
$ws.Range("J10").Value = 0.2249750906731122
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.783357333333333
$ws.Range("N10").Value = 5.350072
$ws.Range("O10").Value = 0.01725253648618087
$ws.Range("P10").Value = 0.01725253648618087
$ws.Range("Q10").Value = 66.65071309171023
$ws.Range("R10").Value = 599.856417825392
$ws.Range("S10").Value = 0.003881390960319718
$ws.Range("T10").Value = 0.003881390960319718
